# Actualización de los planes individuales
# Updates the "Criterio de entrada" (column C) text for several rows of the
# "schedule" sheet, clears two of them, and repositions the sheet's view
# (topLeftCell + active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")
if (-not $ws) { $ws = $wb.ActiveSheet }

$ws.Range("C8").Value = "El equipo ha completado los productos especificados. El equipo ha acumulado toda la información y ha completado todas las formas requeridas."
$ws.Range("C2").Value = "El equipo ha completado un ciclo preeviamente."
$ws.Range("C3").Value = "El equipo ha completado un ciclo preeviamente. Cada miembro del equipo ha leído el capítulo correspondiente a su rol."
$ws.Range("C4").Value = "El equipo ha completado y actualizado el plan de un ciclo preevio."
$ws.Range("C5").Value = "Se ha creado el reporte de inspección del documento de arquitectura."
$ws.Range("C6").ClearContents()
$ws.Range("C7").ClearContents()

$ws.Activate()
$ws.Range("B7").Select()
$excel.ActiveWindow.ScrollColumn = 3
